$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 833
$ws.Range("A2").Value = 50
$ws.Range("A3").Value = 17
$ws.Range("A4").Value = 0.0000703
$ws.Range("A5").Formula = "=67"

$ws.Range("A1:A4").Font.Bold = $false
$ws.Range("A5").Font.Bold = $true

$ws.Range("C8").Select() | Out-Null
